$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.056.32'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.45%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.564.97'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.82%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.90%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.04'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.73%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.491'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.42%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.43%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.08'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.50%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.249'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.50%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0596'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.25%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0860'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.47%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.782.71'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.56%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.564.66'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.79%  '

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.69%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.04%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.041.56'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.49%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.95'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.54%  '

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.37%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.66'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.68%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.40%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.15%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.17'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.24%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.32%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.64'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.06%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.80%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.08'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.77%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.00%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.17%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.78%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0472'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.41%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.94%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.432.17'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.56%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +15.77%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.85%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.34'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.13%  '

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.58%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.532'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.58%  '

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.44'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.00%  '

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.84'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.44%  '

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.808'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.09%  '

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.35%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.15%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.39'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.08%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.56%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.701.26'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.87%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.59'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.08%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.77%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0517'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.36%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0959'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.02%  '
